# Apply updated crypto price/volume values per commit diff.
# Cells store numeric-looking prices as TEXT (matches source file,
# which uses inlineStr for every D/E cell). To keep Excel from
# auto-converting values like "1.002" into numbers, force the cell
# to text format before writing, then restore the default "Normal"
# style so no stray number-format style is left behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue "D2" "26.112.75"
Set-TextValue "D3" "1.666.30"
Set-TextValue "E3" "  -1.26%  "
Set-TextValue "D4" "1.002"
Set-TextValue "E4" "  -0.81%  "
Set-TextValue "D5" "209.76"
Set-TextValue "E5" "  -4.10%  "
Set-TextValue "D6" "0.5183"
Set-TextValue "E6" "  -4.94%  "
Set-TextValue "D7" "1.002"
Set-TextValue "E7" "  -0.80%  "
Set-TextValue "D8" "0.2632"
Set-TextValue "E8" "  -4.11%  "
Set-TextValue "D9" "0.06223"
Set-TextValue "E9" "  -3.47%  "
Set-TextValue "D10" "21.12"
Set-TextValue "E10" "  -3.79%  "
Set-TextValue "D11" "0.07488"
Set-TextValue "E11" "  -2.54%  "
Set-TextValue "D12" "1.769.45"
Set-TextValue "E12" "  +4.88%  "
Set-TextValue "D13" "4.417"
Set-TextValue "E13" "  -2.39%  "
Set-TextValue "D14" "0.5584"
Set-TextValue "E14" "  -4.05%  "
Set-TextValue "D15" "65.94"
Set-TextValue "E15" "  +1.25%  "
Set-TextValue "D16" "0.000007866"
Set-TextValue "E16" "  -6.26%  "
Set-TextValue "D17" "26.139.21"
Set-TextValue "E17" "  -1.02%  "
Set-TextValue "E18" "  -0.77%  "
Set-TextValue "D19" "4.777"
Set-TextValue "E19" "  -3.24%  "
Set-TextValue "D20" "10.36"
Set-TextValue "E20" "  -5.44%  "
Set-TextValue "D21" "186.22"
Set-TextValue "E21" "  -2.69%  "
Set-TextValue "D22" "6.165"
Set-TextValue "E22" "  -1.39%  "
Set-TextValue "D23" "1.003"
Set-TextValue "E23" "  -0.79%  "
Set-TextValue "D24" "147.67"
Set-TextValue "E24" "  -1.24%  "
Set-TextValue "D25" "0.1243"
Set-TextValue "E25" "  -6.08%  "
Set-TextValue "D26" "7.541"
Set-TextValue "E26" "  -4.24%  "
Set-TextValue "D27" "15.91"
Set-TextValue "E27" "  +1.06%  "
Set-TextValue "D28" "0.06242"
Set-TextValue "E28" "  -1.78%  "
Set-TextValue "E29" "  -2.99%  "
Set-TextValue "D30" "1.272"
Set-TextValue "E30" "  -4.20%  "
Set-TextValue "D31" "3.469"
Set-TextValue "E31" "  -2.93%  "
Set-TextValue "D32" "3.418"
Set-TextValue "E32" "  -4.90%  "
Set-TextValue "D33" "1.622"
Set-TextValue "E33" "  -3.60%  "
Set-TextValue "D34" "0.9957"
Set-TextValue "E34" "  -4.35%  "
Set-TextValue "B35" "ImmutableX"
Set-TextValue "C35" "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue "D35" "0.6016"
Set-TextValue "E35" "  -2.18%  "
Set-TextValue "B36" "HuobiToken"
Set-TextValue "C36" "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextValue "D36" "2.406"
Set-TextValue "E36" "  -0.26%  "
Set-TextValue "D37" "2.702"
Set-TextValue "E37" "  -0.25%  "
Set-TextValue "D38" "6.119"
Set-TextValue "E38" "  -2.38%  "
Set-TextValue "E39" "  -1.74%  "
Set-TextValue "D40" "1.074.31"
Set-TextValue "E40" "  -3.88%  "
Set-TextValue "D41" "0.8612"
Set-TextValue "E41" "  -1.68%  "
Set-TextValue "D42" "1.003"
Set-TextValue "E42" "  -1.17%  "
Set-TextValue "D43" "99.05"
Set-TextValue "E43" "  -2.65%  "
Set-TextValue "D44" "1.813.11"
Set-TextValue "E44" "  -1.35%  "
Set-TextValue "D45" "0.00000000108"
Set-TextValue "E45" "  +1.03%  "
Set-TextValue "D46" "55.97"
Set-TextValue "E46" "  -2.56%  "
Set-TextValue "E47" "  -1.35%  "
Set-TextValue "D48" "0.05251"
Set-TextValue "E48" "  -0.40%  "
Set-TextValue "D49" "7.904"
Set-TextValue "E49" "  -3.42%  "
Set-TextValue "D50" "0.4249"
Set-TextValue "E50" "  -1.30%  "
Set-TextValue "D51" "5.920"
Set-TextValue "E51" "  -2.47%  "
